# Update the "修改时间" (last-modified) timestamp column on every portfolio
# sheet from 202509211537 -> 202509211541, keeping the cells as text
# (matching the workbook's existing numberStoredAsText convention).
#
# A plain `.Value = "202509211541"` assignment would make Excel's COM
# layer auto-detect the digit-only string as a number and silently drop
# its text-ness. To avoid that we stage the new text through a
# self-referential text formula (="202509211541"), which is always typed
# as a string result, and then Copy / PasteSpecial(xlPasteValues) the
# range onto itself. That collapses the formula down to a plain value
# while preserving the string type and leaves styles/number formats
# untouched.

$wb = $excel.ActiveWorkbook

$targets = @(
    @{ Sheet = 1; Col = "E"; FirstRow = 2; LastRow = 9 },
    @{ Sheet = 2; Col = "E"; FirstRow = 2; LastRow = 11 },
    @{ Sheet = 3; Col = "G"; FirstRow = 2; LastRow = 13 }
)

$newValue = "202509211541"
$xlPasteValues = -4163

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)

    for ($r = $t.FirstRow; $r -le $t.LastRow; $r++) {
        $addr = $t.Col + $r
        $ws.Range($addr).Formula = "=""" + $newValue + """"
    }

    $rangeAddr = $t.Col + $t.FirstRow + ":" + $t.Col + $t.LastRow
    $range = $ws.Range($rangeAddr)
    $range.Copy()
    $range.PasteSpecial($xlPasteValues)
}

$excel.CutCopyMode = $false
